$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally holds 20 data rows (rows 2-21, columns A:C).
# The edit prepends 7 new rows of data before the old data (pushing the
# old rows down to 9-28) and appends 3 new rows after (rows 29-31).

# 1) Shift the existing data rows 2-21 down to rows 9-28.
#    Walk bottom-up so we never overwrite a source row before reading it.
for ($r = 21; $r -ge 2; $r--) {
  $destRow = $r + 7
  for ($c = 1; $c -le 3; $c++) {
    $ws.Cells.Item($destRow, $c).Value = $ws.Cells.Item($r, $c).Value2
  }
}

# 2) Write the 7 new rows into the now-vacated rows 2-8.
$newTop = @(
  @(0.5757570266723633, -0.6172752380371094, -0.5019410252571106),
  @(0.4425497055053711, -0.6991405487060547, -0.699306845664978),
  @(0.3754444122314453, -0.6968369483947754, -0.6064528226852417),
  @(0.4661340713500976, -0.7378168106079102, -0.8263083696365356),
  @(0.3465394973754883, -0.7457756996154785, -0.6906525492668152),
  @(0.1021490097045898, -0.7542791366577148, -0.5537225604057312),
  @(0.2407388687133789, -0.6835846900939941, -0.2614910900592804)
)
for ($i = 0; $i -lt 7; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 1).Value = $newTop[$i][0]
  $ws.Cells.Item($row, 2).Value = $newTop[$i][1]
  $ws.Cells.Item($row, 3).Value = $newTop[$i][2]
}

# 3) Append the 3 new rows after the shifted data, as rows 29-31.
$newBottom = @(
  @(-0.4794178009033203, -1.04423999786377, -0.2169336676597595),
  @(0.1147146224975586, -0.7852307558059692, -0.0405309796333313),
  @(-0.0877876281738281, -0.7322115302085876, -0.3498360514640808)
)
for ($i = 0; $i -lt 3; $i++) {
  $row = 29 + $i
  $ws.Cells.Item($row, 1).Value = $newBottom[$i][0]
  $ws.Cells.Item($row, 2).Value = $newBottom[$i][1]
  $ws.Cells.Item($row, 3).Value = $newBottom[$i][2]
}
